$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column I: "Expected Return.1" (duplicate header for a new due-date column) ---
# Copy H1's formatting (bold header style w/ border + centered alignment) onto I1, then set its text.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("I1").Value = "Expected Return.1"

# --- Row 2: Huckleberry Finn ---
$ws.Range("F2").Value = "no"
$ws.Range("H2").Value = 45882
$ws.Range("I2").Font.Bold = $false   # materialize empty cell I2 without altering style

# --- Row 3: To Kill a Mockingbird ---
$ws.Range("F3").Value = "no"
$ws.Range("I3").Font.Bold = $false

# --- Row 4: Atomic Habits ---
$ws.Range("F4").Value = "yes"
$ws.Range("I4").Font.Bold = $false

# --- Row 5: 50 Shades of Grey ---
$ws.Range("F5").Value = "yes"
$ws.Range("H5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H5").Value = 45881
$ws.Range("I5").Font.Bold = $false

# --- Row 6: Defining Decade ---
$ws.Range("F6").Value = "yes"
$ws.Range("H6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H6").Value = 45880
$ws.Range("I6").Font.Bold = $false

# --- Row 7: The Great Gatsby ---
$ws.Range("F7").Value = "yes"
$ws.Range("H7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H7").Value = 45881
$ws.Range("I7").Font.Bold = $false

# --- Row 8: One Hundred Years of Solitude ---
$ws.Range("F8").Value = "no"
$ws.Range("G8").Value = "admin"
$ws.Range("H8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H8").Value = 45880
$ws.Range("I8").Font.Bold = $false

# --- Row 9: Wildblood ---
$ws.Range("F9").Value = "No"
$ws.Range("G9").Value = "JohnDoe"
$ws.Range("H9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H9").Value = 45895
$ws.Range("I9").Font.Bold = $false

# --- Row 10: Algebra of Happiness ---
$ws.Range("F10").Value = "yes"
$ws.Range("H10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H10").Value = 45880
$ws.Range("I10").Font.Bold = $false

# --- Row 11: Algebra of Wealth ---
$ws.Range("F11").Value = "yes"
$ws.Range("H11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H11").Value = 45880
$ws.Range("I11").Font.Bold = $false

# --- Row 12: Oh, the Places You'll Go! ---
$ws.Range("F12").Value = "no"
$ws.Range("I12").Font.Bold = $false

# --- Row 13: Who Moved My Cheese? ---
$ws.Range("F13").Value = "yes"
$ws.Range("I13").Font.Bold = $false

# --- Row 14: Everything Is F*cked ---
$ws.Range("F14").Value = "yes"
$ws.Range("I14").Font.Bold = $false

# --- Row 15: #Girlboss ---
$ws.Range("F15").Value = "yes"
$ws.Range("I15").Font.Bold = $false
